# Weekly update: insert 5 new price rows for "Comercializadora del Agro de
# Limarí - Ciruela" at the top of the data block (row 2), pushing the
# existing rows down by 5. Matches commit: "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at row 2 (before the first data row), using
# formatting from the row below (the old row 2) so we don't inherit the
# bold header style from row 1.
$ws.Range("A2:A6").EntireRow.Insert(-4121, 1)

# The insert still carries over the "below" row's style object as a
# whole; strip it back to the workbook default so the new rows start
# unstyled, matching a freshly-appended data row.
$ws.Range("A2:T6").ClearFormats()

# Column D holds dates rendered via a custom numFmt; restore that on the
# new rows (all other columns keep the plain/default style).
$ws.Range("D2:D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 2 ---
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44616
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103002
$ws.Range("J2").Value = "Ciruela"
$ws.Range("K2").Value = "Angeleno"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 220000
$ws.Range("O2").Value = 230000
$ws.Range("P2").Value = 225000
$ws.Range("Q2").Value = "$/bins (450 kilos)"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 500
$ws.Range("T2").Value = 450

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44616
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103002
$ws.Range("J3").Value = "Ciruela"
$ws.Range("K3").Value = "Angeleno"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 16
$ws.Range("N3").Value = 200000
$ws.Range("O3").Value = 210000
$ws.Range("P3").Value = 205000
$ws.Range("Q3").Value = "$/bins (450 kilos)"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 456
$ws.Range("T3").Value = 450

# --- Row 4 ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44616
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103002
$ws.Range("J4").Value = "Ciruela"
$ws.Range("K4").Value = "Angeleno"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 10
$ws.Range("N4").Value = 160000
$ws.Range("O4").Value = 170000
$ws.Range("P4").Value = 165000
$ws.Range("Q4").Value = "$/bins (450 kilos)"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 367
$ws.Range("T4").Value = 450

# --- Row 5 ---
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44616
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100103
$ws.Range("H5").Value = "Frutos de hueso (carozo)"
$ws.Range("I5").Value = 100103002
$ws.Range("J5").Value = "Ciruela"
$ws.Range("K5").Value = "Black Amber"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 200000
$ws.Range("O5").Value = 210000
$ws.Range("P5").Value = 205000
$ws.Range("Q5").Value = "$/bins (450 kilos)"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 456
$ws.Range("T5").Value = 450

# --- Row 6 ---
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44616
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100103
$ws.Range("H6").Value = "Frutos de hueso (carozo)"
$ws.Range("I6").Value = 100103002
$ws.Range("J6").Value = "Ciruela"
$ws.Range("K6").Value = "Black Amber"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 150000
$ws.Range("O6").Value = 160000
$ws.Range("P6").Value = 155000
$ws.Range("Q6").Value = "$/bins (450 kilos)"
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 344
$ws.Range("T6").Value = 450

